# correct some ui things
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update color values in column D (rows 4-6)
$ws.Range("D4").Value = "#ffcbdb"
$ws.Range("D5").Value = "#808080"
$ws.Range("D6").Value = "#52442c"

# Update the selected/active cell to D6
$ws.Range("D6").Select()
